# Daily refresh of the cryptos price list (GitHub Actions scheduled update).
# Updates the Price (D) and Volume(1h) (E) columns for every coin row, and
# reflects the reordering of TheSandbox / InternetComputer(DFINITY) rows.
#
# Price cells that look like plain numbers ("1.002", "315.11", ...) need to
# stay TEXT (matching the source data's inline-string cells) instead of
# Excel's default behaviour of auto-converting such input into a number.
# We do that by momentarily flipping the cell to a Text number format,
# assigning the literal string, then restoring the default "Normal" style
# so the cell's style index is left unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.553.29"
$ws.Range("E2").Value = "  +1.40%  "

# Row 3
$ws.Range("D3").Value = "1.911.49"
$ws.Range("E3").Value = "  +4.72%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.24%  "

# Row 6
$ws.Range("E6").Value = "  -0.07%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5165"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.99%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3956"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.87%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09648"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.84%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.151"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.79%  "

# Row 11
$ws.Range("E11").Value = "  +1.97%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.524"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.18%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.85%  "

# Row 14
$ws.Range("D14").Value = "1.918.44"
$ws.Range("E14").Value = "  +5.56%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.505"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.55%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.01%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.14%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001134"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.54%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06646"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.31%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.20%  "

# Row 21
$ws.Range("E21").Value = "  -0.05%  "

# Row 22
$ws.Range("E22").Value = "  +4.80%  "

# Row 23
$ws.Range("D23").Value = "28.621.06"
$ws.Range("E23").Value = "  +1.46%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.31%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.313"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.97%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.674"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.19%  "

# Row 27
$ws.Range("D27").Value = "2.133.87"
$ws.Range("E27").Value = "  +5.15%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.10%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.40%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.35%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.109"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.53%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1079"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.42%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.752"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.45%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.634"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.49%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.959"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.41%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06799"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.77%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.281"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.86%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02435"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.77%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2220"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.03%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.40%  "

# Row 41
$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.088"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.20%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6462"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.69%  "

# Row 43
$ws.Range("E43").Value = "  +0.31%  "

# Row 44
$ws.Range("E44").Value = "  -0.07%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.29%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6101"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.61%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.777"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.80%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.282"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.50%  "

# Row 49
$ws.Range("E49").Value = "  +4.25%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "125.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.76%  "

# Row 51
$ws.Range("E51").Value = "  +1.60%  "
